$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Update text (shared-string) columns A-H for rows 2-12 with refreshed MeltR method-1 results ----
$ws.Range("A2").Value = "ACCGGU"
$ws.Range("B2").Value = "-54.89 (-57.11 to -51.93)"
$ws.Range("C2").Value = "-54.1 (Â±4)"
$ws.Range("D2").Value = "-154.2 (-161.22 to -144.69)"
$ws.Range("E2").Value = "-151.5 (Â±12.4)"
$ws.Range("F2").Value = "-7.07 (-7.12 to -7.01)"
$ws.Range("G2").Value = "-7.11 (Â±0.2)"
$ws.Range("H2").Value = "45.05 (44.6 to 45.49)"
$ws.Range("A3").Value = "CCAUGG"
$ws.Range("B3").Value = "-60.37 (-61.33 to -59.56)"
$ws.Range("C3").Value = "-60.4 (Â±2)"
$ws.Range("D3").Value = "-172.79 (-175.85 to -170.23)"
$ws.Range("E3").Value = "-173.6 (Â±6.4)"
$ws.Range("F3").Value = "-6.78 (-6.8 to -6.75)"
$ws.Range("G3").Value = "-6.52 (Â±0.05)"
$ws.Range("H3").Value = "42.79 (42.65 to 42.88)"
$ws.Range("A4").Value = "CGAAAGGU/ACCUUUCG"
$ws.Range("B4").Value = "-68.2 (-69.59 to -66.9)"
$ws.Range("C4").Value = "-67.02 (Â±3.025)"
$ws.Range("D4").Value = "-188.37 (-192.64 to -184.42)"
$ws.Range("E4").Value = "-184.51 (Â±9.497)"
$ws.Range("F4").Value = "-9.77 (-9.84 to -9.71)"
$ws.Range("G4").Value = "-9.8 (Â±0.087)"
$ws.Range("H4").Value = "52.48 (52.25 to 52.64)"
$ws.Range("A5").Value = "CGCGCG"
$ws.Range("B5").Value = "-55.77 (-56.62 to -55)"
$ws.Range("C5").Value = "-55.6 (Â±1.3)"
$ws.Range("D5").Value = "-153.57 (-156.2 to -151.15)"
$ws.Range("E5").Value = "-153.9 (Â±4.3)"
$ws.Range("F5").Value = "-8.14 (-8.18 to -8.11)"
$ws.Range("G5").Value = "-7.85 (Â±0.17)"
$ws.Range("H5").Value = "51.34 (51.25 to 51.43)"
$ws.Range("A6").Value = "CGUUGC/GCAACG"
$ws.Range("B6").Value = "-51.61 (-52.98 to -49.85)"
$ws.Range("C6").Value = "-51.06 (Â±3.567)"
$ws.Range("D6").Value = "-143.79 (-148.21 to -138.08)"
$ws.Range("E6").Value = "-141.39 (Â±11.476)"
$ws.Range("F6").Value = "-7.01 (-7.07 to -6.94)"
$ws.Range("G6").Value = "-7.21 (Â±0.087)"
$ws.Range("H6").Value = "39.91 (39.46 to 40.3)"
$ws.Range("A7").Value = "CUGAGUC/GACUCAG"
$ws.Range("B7").Value = "-64.24 (-65.12 to -63.31)"
$ws.Range("C7").Value = "-63.32 (Â±1.928)"
$ws.Range("D7").Value = "-177.85 (-180.62 to -174.9)"
$ws.Range("E7").Value = "-174.65 (Â±6.191)"
$ws.Range("F7").Value = "-9.08 (-9.1 to -9.06)"
$ws.Range("G7").Value = "-9.15 (Â±0.03)"
$ws.Range("H7").Value = "49.82 (49.68 to 49.95)"
$ws.Range("A8").Value = "FAMCGAAAGGU/ACCUUUCGBHQ1"
$ws.Range("B8").Value = "-83.62 (-87.32 to -79.99)"
$ws.Range("C8").Value = "-81.12 (Â±1.583)"
$ws.Range("D8").Value = "-225.31 (-236.37 to -214.42)"
$ws.Range("E8").Value = "-217.61 (Â±5.076)"
$ws.Range("F8").Value = "-13.75 (-14.01 to -13.48)"
$ws.Range("G8").Value = "-13.63 (Â±0.092)"
$ws.Range("H8").Value = "66.28 (66.02 to 66.65)"
$ws.Range("A9").Value = "FAMCGUUGC/GCAACGBHQ1"
$ws.Range("B9").Value = "-58.46 (-60.45 to -56.31)"
$ws.Range("C9").Value = "-61.47 (Â±3.4)"
$ws.Range("D9").Value = "-155.88 (-161.97 to -149.38)"
$ws.Range("E9").Value = "-164.74 (Â±10.499)"
$ws.Range("F9").Value = "-10.12 (-10.29 to -9.96)"
$ws.Range("G9").Value = "-10.38 (Â±0.163)"
$ws.Range("H9").Value = "57.26 (56.88 to 58.17)"
$ws.Range("A10").Value = "FAMCUGAGUC/GACUCAGBHQ1"
$ws.Range("B10").Value = "-75.26 (-77.4 to -73.85)"
$ws.Range("C10").Value = "-74.43 (Â±2.316)"
$ws.Range("D10").Value = "-200.22 (-206.44 to -196.06)"
$ws.Range("E10").Value = "-197.46 (Â±6.709)"
$ws.Range("F10").Value = "-13.17 (-13.4 to -13.03)"
$ws.Range("G10").Value = "-13.19 (Â±0.275)"
$ws.Range("H10").Value = "66.99 (66.86 to 67.19)"
$ws.Range("A11").Value = "GAUAUAUC"
$ws.Range("B11").Value = "-72.95 (-74.09 to -71.66)"
$ws.Range("C11").Value = "-74.2 (Â±4.4)"
$ws.Range("D11").Value = "-217.88 (-221.54 to -213.71)"
$ws.Range("E11").Value = "-221.7 (Â±14.2)"
$ws.Range("F11").Value = "-5.37 (-5.4 to -5.33)"
$ws.Range("G11").Value = "-5.41 (Â±0.06)"
$ws.Range("H11").Value = "35.72 (35.58 to 35.84)"
$ws.Range("A12").Value = "GCAAUUGC"
$ws.Range("B12").Value = "-77.51 (-78.71 to -75.8)"
$ws.Range("C12").Value = "-79.4 (Â±3.8)"
$ws.Range("D12").Value = "-222.11 (-225.79 to -216.83)"
$ws.Range("E12").Value = "-229.8 (Â±11.9)"
$ws.Range("F12").Value = "-8.63 (-8.69 to -8.55)"
$ws.Range("G12").Value = "-8.15 (Â±0.12)"
$ws.Range("H12").Value = "49.27 (49.14 to 49.35)"
$ws.Range("A13").Value = "UAUAUAUA"
$ws.Range("B13").Value = "-62.24 (-66.05 to -59.75)"
$ws.Range("C13").Value = "-63.1 (Â±2.2)"
$ws.Range("D13").Value = "-193.76 (-206.92 to -185.22)"
$ws.Range("E13").Value = "-196.1 (Â±7.4)"
$ws.Range("F13").Value = "-2.14 (-2.31 to -1.87)"
$ws.Range("G13").Value = "-2.27 (Â±0.09)"
$ws.Range("H13").Value = "20.35 (20.08 to 20.51)"

# ---- Update numeric percent-error columns J/K/L/M for rows 2-12 ----
$ws.Range("J2").Value = 1.4496742820442201
$ws.Range("K2").Value = 1.76643768400392
$ws.Range("M2").Value = 0.77390823659480701
$ws.Range("J3").Value = 0.049681212221580097
$ws.Range("K3").Value = 0.467680937671412
$ws.Range("M3").Value = 3.3020548758760002
$ws.Range("J4").Value = 1.74530394911996
$ws.Range("K4").Value = 2.0703711649860601
$ws.Range("M4").Value = 0.79711520212564402
$ws.Range("J5").Value = 0.30528867738170401
$ws.Range("K5").Value = 0.214655088301306
$ws.Range("M5").Value = 3.4475926292847299
$ws.Range("J6").Value = 1.0713937859160401
$ws.Range("K6").Value = 1.6831474857984501
$ws.Range("M6").Value = 3.1808654913081198
$ws.Range("J7").Value = 1.44245845092505
$ws.Range("K7").Value = 1.8156028368794299
$ws.Range("M7").Value = 1.15745360207543
$ws.Range("J8").Value = 3.0350855894136202
$ws.Range("K8").Value = 3.47692585568499
$ws.Range("L8").Value = 0.87655222790357301
$ws.Range("J9").Value = 5.0195947636121003
$ws.Range("K9").Value = 5.52679184080844
$ws.Range("L9").Value = 2.5365853658536701
$ws.Range("M9").Value = 0.76548364648574296
$ws.Range("J10").Value = 1.10895851426281
$ws.Range("K10").Value = 1.3880506940253401
$ws.Range("L10").Value = 0.15174506828527801
$ws.Range("M10").Value = 0.75842070042383103
$ws.Range("J11").Value = 1.6989466530750901
$ws.Range("K11").Value = 1.73802265799172
$ws.Range("M11").Value = 0.50265289025411797
$ws.Range("J12").Value = 2.4090242814352201
$ws.Range("K12").Value = 3.4033325219623398
$ws.Range("M12").Value = 4.7159031889477596

# ---- Old row 13 (averages) becomes a normal data row for the new UAUAUAUA helix; drop its old number format ----
$ws.Range("J13:M13").ClearFormats()
$ws.Range("I13").Value = 21.1
$ws.Range("J13").Value = 1.3722674325833699
$ws.Range("K13").Value = 1.2004309239214099
$ws.Range("L13").Value = 5.8956916099773196
$ws.Range("M13").Value = 3.6188178528347401

$ws.Range("J14").Formula = "=AVERAGE(J2:J13)"
$ws.Range("K14:M14").Formula = "=AVERAGE(K2:K13)"

# ---- Selection moved to the new trailing average row ----
$ws.Range("J14:M14").Select() | Out-Null
